$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 (file 64799798-...) gets a freshly generated handback report,
# so its handoff/handback datetimes are updated to the new run's timestamps.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 04:53:04"
$wsZhCn.Range("H2").Value = "2016-03-23 04:53:28"

# de-de sheet: same file, same new handback run, updated timestamps.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 04:53:09"
$wsDeDe.Range("H2").Value = "2016-03-23 04:53:34"
